# Auto-generated Excel COM-interop script to apply the scheduled-runner market-data update
# across the Exodus_Profits workbook's 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H29").Value = 1529.8334
$ws.Range("I29").Value = 44.75
$ws.Range("K29").Value = 134.25
$ws.Range("M29").Value = 146.75
$ws.Range("H62").Value = 5278.9414
$ws.Range("I62").Value = 4334.7144
$ws.Range("J62").Value = 5939.9
$ws.Range("K62").Value = 4334.7144
$ws.Range("L62").Value = 5939.9
$ws.Range("M62").Value = -3710.7144
$ws.Range("N62").Value = -7187.9
$ws.Range("H65").Value = 5278.9414
$ws.Range("I65").Value = 4334.7144
$ws.Range("J65").Value = 5939.9
$ws.Range("K65").Value = 21673.572
$ws.Range("L65").Value = 29699.5
$ws.Range("M65").Value = -18553.572
$ws.Range("N65").Value = -35939.5
$ws.Range("H70").Value = 977.3913
$ws.Range("I70").Value = 899.0769
$ws.Range("K70").Value = 2697.2307
$ws.Range("M70").Value = -2427.2307
$ws.Range("H73").Value = 977.3913
$ws.Range("I73").Value = 899.0769
$ws.Range("K73").Value = 2697.2307
$ws.Range("M73").Value = -1761.2307
$ws.Range("H100").Value = 1491.5
$ws.Range("I100").Value = 770
$ws.Range("K100").Value = 770
$ws.Range("M100").Value = -229
$ws.Range("H125").Value = 1154.4375
$ws.Range("I125").Value = 1143.75
$ws.Range("J125").Value = 1165.125
$ws.Range("K125").Value = 10293.75
$ws.Range("L125").Value = 10486.125
$ws.Range("M125").Value = -7833.75
$ws.Range("N125").Value = -15406.125
$ws.Range("H133").Value = 94009
$ws.Range("J133").Value = 94009
$ws.Range("L133").Value = 94009
$ws.Range("N133").Value = -104129
$ws.Range("H134").Value = 81965.75
$ws.Range("J134").Value = 94825.164
$ws.Range("L134").Value = 94825.164
$ws.Range("N134").Value = -104965.164
$ws.Range("H136").Value = 68169.8
$ws.Range("J136").Value = 81956.86
$ws.Range("L136").Value = 81956.86
$ws.Range("N136").Value = -92156.86
$ws.Range("H139").Value = 99990
$ws.Range("J139").Value = 99990
$ws.Range("L139").Value = 99990
$ws.Range("N139").Value = -110270
$ws.Range("H141").Value = 4597.6
$ws.Range("I141").Value = 4597.6
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 13792.8
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -8612.800000000001
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H32").Value = 9352.440000000001
$ws.Range("I32").Value = 5744.487
$ws.Range("K32").Value = 5744.487
$ws.Range("M32").Value = -5457.487
$ws.Range("H61").Value = 2699.375
$ws.Range("I61").Value = 1207.7
$ws.Range("K61").Value = 1207.7
$ws.Range("M61").Value = -995.7
$ws.Range("H74").Value = 2470.2188
$ws.Range("I74").Value = 1742.5172
$ws.Range("J74").Value = 9504.666999999999
$ws.Range("K74").Value = 1742.5172
$ws.Range("L74").Value = 9504.666999999999
$ws.Range("M74").Value = -868.5172
$ws.Range("N74").Value = -11252.667
$ws.Range("H77").Value = 2470.2188
$ws.Range("I77").Value = 1742.5172
$ws.Range("J77").Value = 9504.666999999999
$ws.Range("K77").Value = 8712.585999999999
$ws.Range("L77").Value = 47523.335
$ws.Range("M77").Value = -4344.585999999999
$ws.Range("N77").Value = -56259.335
$ws.Range("H110").Value = 1394.8572
$ws.Range("I110").Value = 1226.0588
$ws.Range("J110").Value = 2112.25
$ws.Range("K110").Value = 1226.0588
$ws.Range("L110").Value = 2112.25
$ws.Range("M110").Value = 818.9412
$ws.Range("N110").Value = -6202.25
$ws.Range("H122").Value = 6336.8237
$ws.Range("I122").Value = 6485.5
$ws.Range("K122").Value = 19456.5
$ws.Range("M122").Value = -17006.5
$ws.Range("H131").Value = 82499.5
$ws.Range("J131").Value = 82499.5
$ws.Range("L131").Value = 82499.5
$ws.Range("N131").Value = -92579.5
$ws.Range("H132").Value = 3162.0557
$ws.Range("I132").Value = 2869.0715
$ws.Range("J132").Value = 4187.5
$ws.Range("K132").Value = 8607.2145
$ws.Range("L132").Value = 12562.5
$ws.Range("M132").Value = -6077.2145
$ws.Range("N132").Value = -17622.5
$ws.Range("H136").Value = 2699.375
$ws.Range("I136").Value = 1207.7
$ws.Range("K136").Value = 3623.1
$ws.Range("M136").Value = -1073.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 54670.8
$ws.Range("J109").Value = 54670.8
$ws.Range("L109").Value = 54670.8
$ws.Range("N109").Value = -57444.8
$ws.Range("H135").Value = 98392
$ws.Range("J135").Value = 98392
$ws.Range("L135").Value = 98392
$ws.Range("N135").Value = -108532

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 64999.668
$ws.Range("J50").Value = 64999.668
$ws.Range("L50").Value = 64999.668
$ws.Range("N50").Value = -66249.66800000001
$ws.Range("H59").Value = 91034.5
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 91034.5
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 91034.5
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -93324.5
$ws.Range("H60").Value = 24575.166
$ws.Range("I60").Value = 10000
$ws.Range("K60").Value = 10000
$ws.Range("M60").Value = -9489
$ws.Range("H62").Value = 3439.625
$ws.Range("I62").Value = 3073.8572
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 3073.8572
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -2449.8572
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 3439.625
$ws.Range("I65").Value = 3073.8572
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 15369.286
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -12249.286
$ws.Range("N65").Value = -36240
$ws.Range("H103").Value = 9577
$ws.Range("J103").Value = 9577
$ws.Range("L103").Value = 9577
$ws.Range("N103").Value = -11921
$ws.Range("H105").Value = 63891.11
$ws.Range("I105").Value = 112905.2
$ws.Range("J105").Value = 2623.5
$ws.Range("K105").Value = 112905.2
$ws.Range("L105").Value = 2623.5
$ws.Range("M105").Value = -111158.2
$ws.Range("N105").Value = -6117.5
$ws.Range("H132").Value = 2739.1428
$ws.Range("J132").Value = 2639.8
$ws.Range("L132").Value = 7919.400000000001
$ws.Range("N132").Value = -12979.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 17500.5
$ws.Range("J80").Value = 33000
$ws.Range("L80").Value = 99000
$ws.Range("N80").Value = -100872
$ws.Range("H83").Value = 17500.5
$ws.Range("J83").Value = 33000
$ws.Range("L83").Value = 297000
$ws.Range("N83").Value = -306360
$ws.Range("H122").Value = 723363.0600000001
$ws.Range("I122").Value = 1033.3334
$ws.Range("J122").Value = 920362.0600000001
$ws.Range("K122").Value = 9300.000599999999
$ws.Range("L122").Value = 8283258.540000001
$ws.Range("M122").Value = -6850.000599999999
$ws.Range("N122").Value = -8288158.540000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 65182.223
$ws.Range("J108").Value = 65182.223
$ws.Range("L108").Value = 65182.223
$ws.Range("N108").Value = -72862.223
$ws.Range("H132").Value = 2608.6
$ws.Range("I132").Value = 2057.2727
$ws.Range("J132").Value = 4124.75
$ws.Range("K132").Value = 6171.8181
$ws.Range("L132").Value = 12374.25
$ws.Range("M132").Value = -3641.8181
$ws.Range("N132").Value = -17434.25
$ws.Range("H135").Value = 68993
$ws.Range("J135").Value = 68993
$ws.Range("L135").Value = 68993
$ws.Range("N135").Value = -79133
$ws.Range("H140").Value = 44448.11
$ws.Range("J140").Value = 51578.25
$ws.Range("L140").Value = 51578.25
$ws.Range("N140").Value = -61938.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 200
$ws.Range("I53").Value = 200
$ws.Range("J53").Value = 0
$ws.Range("K53").Value = 200
$ws.Range("L53").Value = 0
$ws.Range("M53").Value = 318
$ws.Range("N53").ClearContents()
$ws.Range("H82").Value = 1330.7222
$ws.Range("I82").Value = 717.8570999999999
$ws.Range("J82").Value = 3475.75
$ws.Range("K82").Value = 717.8570999999999
$ws.Range("L82").Value = 3475.75
$ws.Range("M82").Value = -356.8570999999999
$ws.Range("N82").Value = -4197.75
$ws.Range("H85").Value = 1330.7222
$ws.Range("I85").Value = 717.8570999999999
$ws.Range("J85").Value = 3475.75
$ws.Range("K85").Value = 717.8570999999999
$ws.Range("L85").Value = 3475.75
$ws.Range("M85").Value = 530.1429000000001
$ws.Range("N85").Value = -5971.75
$ws.Range("H97").Value = 15948.333
$ws.Range("J97").Value = 15948.333
$ws.Range("L97").Value = 15948.333
$ws.Range("N97").Value = -17930.333
$ws.Range("H117").Value = 62752
$ws.Range("J117").Value = 62752
$ws.Range("L117").Value = 62752
$ws.Range("N117").Value = -71930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1832.9531
$ws.Range("I136").Value = 1286.4
$ws.Range("K136").Value = 3859.2
$ws.Range("M136").Value = -1309.2
